# Add a new "test" worksheet (mirrors the other scenario sheets) with a
# couple of new test rows, make it the active tab, and tidy up the
# selections/tabSelected flags on the sheets that used to hold them.

$wb = $excel.ActiveWorkbook

# The "batch-file" sheet currently carries tabSelected="1"; that marker will
# move to the new "test" sheet once it becomes active below (its own stored
# selection of C14 is left untouched).

# The "logon" sheet's lingering selection (T14) becomes a plain range
# selection over A1:C1.
$logonSheet = $wb.Worksheets.Item("logon")
$logonSheet.Activate()
$logonSheet.Range("A1:C1").Select()

# Create the new "test" sheet at the end of the workbook.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "test"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch the sheet by name: after Move() the old reference no longer
# tracks the sheet's new position, so look it up fresh before editing it.
$testSheet = $wb.Worksheets.Item("test")

# Header row, matching the layout used by every other sheet in the workbook
# (header cells use the shared "Text" number format already present in the
# workbook's style table).
$testSheet.Range("A1").Value = "Scenario No."
$testSheet.Range("B1").Value = "key"
$testSheet.Range("C1").Value = "value"
$testSheet.Range("A1:C1").NumberFormat = "@"

# Data row (note: "c d" must land in the shared-string table before
# "test1" so the indices come out in the same order as the original edit).
$testSheet.Range("A2").Value = "scenario1"
$testSheet.Range("C2").Value = "c d"
$testSheet.Range("B2").Value = "test1"

$testSheet.Activate()
$testSheet.Range("A1:C1").Select()
